$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds serial-date values (style uses date format).
# Rows 2-17 all change from 45184 (2023-09-15) to 45185 (2023-09-16).
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
